$d = $word.ActiveDocument

# 1. Merge "4. " + "Escreva  um" + " " => "4. Escreva  um "
$d.Content.Find.Execute("4. Escreva  um", $true, $false, $false, $false, $false, $true, 1, $false, "4. Escreva  um", 2)

# 2. Split "2,3 e 5" into "2," + " " + "3 e 5" (insert space)
$d.Content.Find.Execute("o peso das notas é: 2,3 e 5, respectivamente.", $true, $false, $false, $false, $false, $true, 1, $false, "o peso das notas é: 2, 3 e 5, respectivamente.", 2)

# 3. Merge ", pode ser resolvido segundo mostrado " + "abaixo :" + " " => single run
$d.Content.Find.Execute(", pode ser resolvido segundo mostrado abaixo : ", $true, $false, $false, $false, $false, $true, 1, $false, ", pode ser resolvido segundo mostrado abaixo : ", 2)

# 4. Merge " que lê os coeficientes " + "a,b" + ",c,d,e e f e calcula e mostra os valores de x e y. " => single run
$d.Content.Find.Execute(" que lê os coeficientes a,b,c,d,e e f e calcula e mostra os valores de x e y. ", $true, $false, $false, $false, $false, $true, 1, $false, " que lê os coeficientes a,b,c,d,e e f e calcula e mostra os valores de x e y. ", 2)
